$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting old D:K data to F:M
$ws.Range("D:E").Insert()

# Copy number formatting from the (now-shifted) F:G columns, which retain the original
# D:E formatting, onto the newly inserted D:E columns across the data blocks only
# (avoiding blank separator rows 36-37 and 78-79 so no stray cells are introduced there).
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate cell values for quarterly financial data (columns D through M)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43282
$ws.Range("G7").Value = 43191
$ws.Range("H7").Value = 43100
$ws.Range("I7").Value = 43009
$ws.Range("J7").Value = 42918
$ws.Range("K7").Value = 42827
$ws.Range("L7").Value = 42735
$ws.Range("M7").Value = 42645
$ws.Range("D8").Value = 1987900
$ws.Range("E8").Value = 2079600
$ws.Range("F8").Value = 1751600
$ws.Range("G8").Value = 1972000
$ws.Range("H8").Value = 1939600
$ws.Range("I8").Value = 2033100
$ws.Range("J8").Value = 1663000
$ws.Range("K8").Value = 1879700
$ws.Range("L8").Value = 1970200
$ws.Range("M8").Value = 2003500
$ws.Range("D9").Value = 1046300
$ws.Range("E9").Value = 1211500
$ws.Range("F9").Value = 950900
$ws.Range("G9").Value = 995700
$ws.Range("H9").Value = 1103700
$ws.Range("I9").Value = 1090000
$ws.Range("J9").Value = 891400
$ws.Range("K9").Value = 970300
$ws.Range("L9").Value = 1228000
$ws.Range("M9").Value = 1152600
$ws.Range("D10").Value = 941600
$ws.Range("E10").Value = 868100
$ws.Range("F10").Value = 800700
$ws.Range("G10").Value = 976300
$ws.Range("H10").Value = 835900
$ws.Range("I10").Value = 943100
$ws.Range("J10").Value = 771600
$ws.Range("K10").Value = 909400
$ws.Range("L10").Value = 742200
$ws.Range("M10").Value = 850900
$ws.Range("D11").Value = $null
$ws.Range("E11").Value = $null
$ws.Range("F11").Value = $null
$ws.Range("G11").Value = $null
$ws.Range("H11").Value = $null
$ws.Range("I11").Value = $null
$ws.Range("J11").Value = $null
$ws.Range("K11").Value = $null
$ws.Range("L11").Value = $null
$ws.Range("M11").Value = $null
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("L12").Value = "NA"
$ws.Range("M12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("D14").Value = 65900
$ws.Range("E14").Value = 15800
$ws.Range("F14").Value = 61800
$ws.Range("G14").Value = 16400
$ws.Range("H14").Value = 41900
$ws.Range("I14").Value = 22000
$ws.Range("J14").Value = 24700
$ws.Range("K14").Value = 252700
$ws.Range("L14").Value = 14800
$ws.Range("M14").Value = 24300
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("D16").Value = $null
$ws.Range("E16").Value = $null
$ws.Range("F16").Value = $null
$ws.Range("G16").Value = $null
$ws.Range("H16").Value = $null
$ws.Range("I16").Value = $null
$ws.Range("J16").Value = $null
$ws.Range("K16").Value = $null
$ws.Range("L16").Value = $null
$ws.Range("M16").Value = $null
$ws.Range("D17").Value = 1594000
$ws.Range("E17").Value = 1676600
$ws.Range("F17").Value = 1455200
$ws.Range("G17").Value = 1491900
$ws.Range("H17").Value = 1647600
$ws.Range("I17").Value = 1585400
$ws.Range("J17").Value = 1352800
$ws.Range("K17").Value = 1682400
$ws.Range("L17").Value = 1748400
$ws.Range("M17").Value = 1650200
$ws.Range("D18").Value = 393900
$ws.Range("E18").Value = 403000
$ws.Range("F18").Value = 296400
$ws.Range("G18").Value = 480100
$ws.Range("H18").Value = 292000
$ws.Range("I18").Value = 447700
$ws.Range("J18").Value = 310200
$ws.Range("K18").Value = 197300
$ws.Range("L18").Value = 221800
$ws.Range("M18").Value = 353300
$ws.Range("D19").Value = $null
$ws.Range("E19").Value = $null
$ws.Range("F19").Value = $null
$ws.Range("G19").Value = $null
$ws.Range("H19").Value = $null
$ws.Range("I19").Value = $null
$ws.Range("J19").Value = $null
$ws.Range("K19").Value = $null
$ws.Range("L19").Value = $null
$ws.Range("M19").Value = $null
$ws.Range("D20").Value = -10100
$ws.Range("E20").Value = -7400
$ws.Range("F20").Value = 400
$ws.Range("G20").Value = 700
$ws.Range("H20").Value = -4800
$ws.Range("I20").Value = -22000
$ws.Range("J20").Value = -4700
$ws.Range("K20").Value = -4900
$ws.Range("L20").Value = 1000
$ws.Range("M20").Value = -700
$ws.Range("D21").Value = 453100
$ws.Range("E21").Value = 470400
$ws.Range("F21").Value = 373400
$ws.Range("G21").Value = 555200
$ws.Range("H21").Value = 354700
$ws.Range("I21").Value = 488000
$ws.Range("J21").Value = 372700
$ws.Range("K21").Value = 257300
$ws.Range("L21").Value = 282800
$ws.Range("M21").Value = 437700
$ws.Range("D22").Value = 39800
$ws.Range("E22").Value = 38700
$ws.Range("F22").Value = 36800
$ws.Range("G22").Value = 31600
$ws.Range("H22").Value = 26800
$ws.Range("I22").Value = 24900
$ws.Range("J22").Value = 24400
$ws.Range("K22").Value = 24000
$ws.Range("L22").Value = 24200
$ws.Range("M22").Value = 24700
$ws.Range("D23").Value = 344000
$ws.Range("E23").Value = 356900
$ws.Range("F23").Value = 260000
$ws.Range("G23").Value = 449200
$ws.Range("H23").Value = 260400
$ws.Range("I23").Value = 400800
$ws.Range("J23").Value = 281100
$ws.Range("K23").Value = 168400
$ws.Range("L23").Value = 198600
$ws.Range("M23").Value = 327800
$ws.Range("D24").Value = 20200
$ws.Range("E24").Value = 91400
$ws.Range("F24").Value = 36700
$ws.Range("G24").Value = 98500
$ws.Range("H24").Value = 78900
$ws.Range("I24").Value = 126800
$ws.Range("J24").Value = 78400
$ws.Range("K24").Value = 70100
$ws.Range("L24").Value = 81800
$ws.Range("M24").Value = 100400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("D26").Value = 323800
$ws.Range("E26").Value = 265400
$ws.Range("F26").Value = 223300
$ws.Range("G26").Value = 350700
$ws.Range("H26").Value = 181500
$ws.Range("I26").Value = 274000
$ws.Range("J26").Value = 202700
$ws.Range("K26").Value = 98200
$ws.Range("L26").Value = 116900
$ws.Range("M26").Value = 227400
$ws.Range("D27").Value = 329000
$ws.Range("E27").Value = 263700
$ws.Range("F27").Value = 226900
$ws.Range("G27").Value = 350200
$ws.Range("H27").Value = 181100
$ws.Range("I27").Value = 273300
$ws.Range("J27").Value = 203500
$ws.Range("K27").Value = 125000
$ws.Range("L27").Value = 116900
$ws.Range("M27").Value = 227400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("D29").Value = 7800
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"
$ws.Range("K29").Value = "NA"
$ws.Range("L29").Value = "NA"
$ws.Range("M29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("D32").Value = 10100
$ws.Range("E32").Value = 7400
$ws.Range("F32").Value = -400
$ws.Range("G32").Value = -700
$ws.Range("H32").Value = 4800
$ws.Range("I32").Value = 22000
$ws.Range("J32").Value = 4700
$ws.Range("K32").Value = 4900
$ws.Range("L32").Value = -1000
$ws.Range("M32").Value = 700
$ws.Range("D33").Value = 336800
$ws.Range("E33").Value = 263700
$ws.Range("F33").Value = 226900
$ws.Range("G33").Value = 350200
$ws.Range("H33").Value = 181100
$ws.Range("I33").Value = 273300
$ws.Range("J33").Value = 203500
$ws.Range("K33").Value = 125000
$ws.Range("L33").Value = 116900
$ws.Range("M33").Value = 227400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("D35").Value = 336800
$ws.Range("E35").Value = 263700
$ws.Range("F35").Value = 226900
$ws.Range("G35").Value = 350200
$ws.Range("H35").Value = 181100
$ws.Range("I35").Value = 273300
$ws.Range("J35").Value = 203500
$ws.Range("K35").Value = 125000
$ws.Range("L35").Value = 116900
$ws.Range("M35").Value = 227400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43282
$ws.Range("G38").Value = 43191
$ws.Range("H38").Value = 43100
$ws.Range("I38").Value = 43009
$ws.Range("J38").Value = 42918
$ws.Range("K38").Value = 42827
$ws.Range("L38").Value = 42735
$ws.Range("M38").Value = 42645
$ws.Range("D39").Value = $null
$ws.Range("E39").Value = $null
$ws.Range("F39").Value = $null
$ws.Range("G39").Value = $null
$ws.Range("H39").Value = $null
$ws.Range("I39").Value = $null
$ws.Range("J39").Value = $null
$ws.Range("K39").Value = $null
$ws.Range("L39").Value = $null
$ws.Range("M39").Value = $null
$ws.Range("D40").Value = $null
$ws.Range("E40").Value = $null
$ws.Range("F40").Value = $null
$ws.Range("G40").Value = $null
$ws.Range("H40").Value = $null
$ws.Range("I40").Value = $null
$ws.Range("J40").Value = $null
$ws.Range("K40").Value = $null
$ws.Range("L40").Value = $null
$ws.Range("M40").Value = $null
$ws.Range("D41").Value = 588000
$ws.Range("E41").Value = 823800
$ws.Range("F41").Value = 467400
$ws.Range("G41").Value = 476400
$ws.Range("H41").Value = 380200
$ws.Range("I41").Value = 275100
$ws.Range("J41").Value = 214100
$ws.Range("K41").Value = 235000
$ws.Range("L41").Value = 297000
$ws.Range("M41").Value = 333300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("D43").Value = 594100
$ws.Range("E43").Value = 814900
$ws.Range("F43").Value = 501900
$ws.Range("G43").Value = 614300
$ws.Range("H43").Value = 588300
$ws.Range("I43").Value = 742800
$ws.Range("J43").Value = 417500
$ws.Range("K43").Value = 595800
$ws.Range("L43").Value = 581400
$ws.Range("M43").Value = 759600
$ws.Range("D44").Value = 784900
$ws.Range("E44").Value = 880700
$ws.Range("F44").Value = 916400
$ws.Range("G44").Value = 782500
$ws.Range("H44").Value = 752800
$ws.Range("I44").Value = 938200
$ws.Range("J44").Value = 936400
$ws.Range("K44").Value = 795400
$ws.Range("L44").Value = 745700
$ws.Range("M44").Value = 843500
$ws.Range("D45").Value = 272200
$ws.Range("E45").Value = 274700
$ws.Range("F45").Value = 479500
$ws.Range("G45").Value = 397300
$ws.Range("H45").Value = 280600
$ws.Range("I45").Value = 258400
$ws.Range("J45").Value = 343600
$ws.Range("K45").Value = 247600
$ws.Range("L45").Value = 192800
$ws.Range("M45").Value = 194000
$ws.Range("D46").Value = 2239200
$ws.Range("E46").Value = 2794000
$ws.Range("F46").Value = 2365100
$ws.Range("G46").Value = 2270500
$ws.Range("H46").Value = 2001900
$ws.Range("I46").Value = 2214500
$ws.Range("J46").Value = 1911500
$ws.Range("K46").Value = 1873900
$ws.Range("L46").Value = 1816800
$ws.Range("M46").Value = 2130500
$ws.Range("D47").Value = "NA"
$ws.Range("E47").Value = "NA"
$ws.Range("F47").Value = "NA"
$ws.Range("G47").Value = "NA"
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 1500
$ws.Range("L47").Value = "NA"
$ws.Range("M47").Value = "NA"
$ws.Range("D48").Value = 2130300
$ws.Range("E48").Value = 2092900
$ws.Range("F48").Value = 2083800
$ws.Range("G48").Value = 2119000
$ws.Range("H48").Value = 2106700
$ws.Range("I48").Value = 2050100
$ws.Range("J48").Value = 2033800
$ws.Range("K48").Value = 2050400
$ws.Range("L48").Value = 2177200
$ws.Range("M48").Value = 2159600
$ws.Range("D49").Value = 3205800
$ws.Range("E49").Value = 2806400
$ws.Range("F49").Value = 2806200
$ws.Range("G49").Value = 2786300
$ws.Range("H49").Value = 1295100
$ws.Range("I49").Value = 1297100
$ws.Range("J49").Value = 1297200
$ws.Range("K49").Value = 1292300
$ws.Range("L49").Value = 1305100
$ws.Range("M49").Value = 1326400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("D52").Value = 127800
$ws.Range("E52").Value = 142800
$ws.Range("F52").Value = 150500
$ws.Range("G52").Value = 157000
$ws.Range("H52").Value = 150000
$ws.Range("I52").Value = 93800
$ws.Range("J52").Value = 137700
$ws.Range("K52").Value = 124300
$ws.Range("L52").Value = 225200
$ws.Range("M52").Value = 228900
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("D54").Value = 7703000
$ws.Range("E54").Value = 7836200
$ws.Range("F54").Value = 7405700
$ws.Range("G54").Value = 7332800
$ws.Range("H54").Value = 5553700
$ws.Range("I54").Value = 5655500
$ws.Range("J54").Value = 5380200
$ws.Range("K54").Value = 5342400
$ws.Range("L54").Value = 5524300
$ws.Range("M54").Value = 5845400
$ws.Range("D55").Value = $null
$ws.Range("E55").Value = $null
$ws.Range("F55").Value = $null
$ws.Range("G55").Value = $null
$ws.Range("H55").Value = $null
$ws.Range("I55").Value = $null
$ws.Range("J55").Value = $null
$ws.Range("K55").Value = $null
$ws.Range("L55").Value = $null
$ws.Range("M55").Value = $null
$ws.Range("D56").Value = $null
$ws.Range("E56").Value = $null
$ws.Range("F56").Value = $null
$ws.Range("G56").Value = $null
$ws.Range("H56").Value = $null
$ws.Range("I56").Value = $null
$ws.Range("J56").Value = $null
$ws.Range("K56").Value = $null
$ws.Range("L56").Value = $null
$ws.Range("M56").Value = $null
$ws.Range("D57").Value = 502300
$ws.Range("E57").Value = 581500
$ws.Range("F57").Value = 496900
$ws.Range("G57").Value = 520000
$ws.Range("H57").Value = 523200
$ws.Range("I57").Value = 529400
$ws.Range("J57").Value = 471500
$ws.Range("K57").Value = 499100
$ws.Range("L57").Value = 522500
$ws.Range("M57").Value = 458000
$ws.Range("D58").Value = 1203300
$ws.Range("E58").Value = 1452300
$ws.Range("F58").Value = 1351300
$ws.Range("G58").Value = 2549500
$ws.Range("H58").Value = 859500
$ws.Range("I58").Value = 1115700
$ws.Range("J58").Value = 622100
$ws.Range("K58").Value = 487600
$ws.Range("L58").Value = 632700
$ws.Range("M58").Value = 862400
$ws.Range("D59").Value = 712900
$ws.Range("E59").Value = 725000
$ws.Range("F59").Value = 666100
$ws.Range("G59").Value = 636000
$ws.Range("H59").Value = 693900
$ws.Range("I59").Value = 692500
$ws.Range("J59").Value = 648600
$ws.Range("K59").Value = 724100
$ws.Range("L59").Value = 754200
$ws.Range("M59").Value = 696600
$ws.Range("D60").Value = 2418600
$ws.Range("E60").Value = 2758800
$ws.Range("F60").Value = 2514300
$ws.Range("G60").Value = 3705500
$ws.Range("H60").Value = 2076500
$ws.Range("I60").Value = 2337700
$ws.Range("J60").Value = 1742200
$ws.Range("K60").Value = 1710900
$ws.Range("L60").Value = 1909400
$ws.Range("M60").Value = 2017000
$ws.Range("D61").Value = 3254300
$ws.Range("E61").Value = 3253900
$ws.Range("F61").Value = 3249700
$ws.Range("G61").Value = 2059900
$ws.Range("H61").Value = 2061000
$ws.Range("I61").Value = 2054100
$ws.Range("J61").Value = 2349800
$ws.Range("K61").Value = 2350900
$ws.Range("L61").Value = 2347500
$ws.Range("M61").Value = 2362500
$ws.Range("D62").Value = 622900
$ws.Range("E62").Value = 574200
$ws.Range("F62").Value = 577900
$ws.Range("G62").Value = 577700
$ws.Range("H62").Value = 484600
$ws.Range("I62").Value = 424700
$ws.Range("J62").Value = 418300
$ws.Range("K62").Value = 435100
$ws.Range("L62").Value = 439700
$ws.Range("M62").Value = 523800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("D66").Value = 6304300
$ws.Range("E66").Value = 6601300
$ws.Range("F66").Value = 6355700
$ws.Range("G66").Value = 6360700
$ws.Range("H66").Value = 4638400
$ws.Range("I66").Value = 4832200
$ws.Range("J66").Value = 4524900
$ws.Range("K66").Value = 4512300
$ws.Range("L66").Value = 4738500
$ws.Range("M66").Value = 4950000
$ws.Range("D67").Value = $null
$ws.Range("E67").Value = $null
$ws.Range("F67").Value = $null
$ws.Range("G67").Value = $null
$ws.Range("H67").Value = $null
$ws.Range("I67").Value = $null
$ws.Range("J67").Value = $null
$ws.Range("K67").Value = $null
$ws.Range("L67").Value = $null
$ws.Range("M67").Value = $null
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("D72").Value = 7032000
$ws.Range("E72").Value = 6843100
$ws.Range("F72").Value = 6727100
$ws.Range("G72").Value = 6634300
$ws.Range("H72").Value = 6371100
$ws.Range("I72").Value = 6325000
$ws.Range("J72").Value = 6187400
$ws.Range("K72").Value = 6112500
$ws.Range("L72").Value = 6116000
$ws.Range("M72").Value = 6129100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("D76").Value = 1398700
$ws.Range("E76").Value = 1234900
$ws.Range("F76").Value = 1049900
$ws.Range("G76").Value = 972100
$ws.Range("H76").Value = 915300
$ws.Range("I76").Value = 823300
$ws.Range("J76").Value = 855300
$ws.Range("K76").Value = 830100
$ws.Range("L76").Value = 785900
$ws.Range("M76").Value = 895400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43282
$ws.Range("G80").Value = 43191
$ws.Range("H80").Value = 43100
$ws.Range("I80").Value = 43009
$ws.Range("J80").Value = 42918
$ws.Range("K80").Value = 42827
$ws.Range("L80").Value = 42735
$ws.Range("M80").Value = 42645
$ws.Range("D81").Value = 336800
$ws.Range("E81").Value = 263700
$ws.Range("F81").Value = 226900
$ws.Range("G81").Value = 350200
$ws.Range("H81").Value = 181100
$ws.Range("I81").Value = 273300
$ws.Range("J81").Value = 203500
$ws.Range("K81").Value = 125000
$ws.Range("L81").Value = 116900
$ws.Range("M81").Value = 227400
$ws.Range("D82").Value = $null
$ws.Range("E82").Value = $null
$ws.Range("F82").Value = $null
$ws.Range("G82").Value = $null
$ws.Range("H82").Value = $null
$ws.Range("I82").Value = $null
$ws.Range("J82").Value = $null
$ws.Range("K82").Value = $null
$ws.Range("L82").Value = $null
$ws.Range("M82").Value = $null
$ws.Range("D83").Value = 69300
$ws.Range("E83").Value = 74800
$ws.Range("F83").Value = 76600
$ws.Range("G83").Value = 74400
$ws.Range("H83").Value = 67500
$ws.Range("I83").Value = 62200
$ws.Range("J83").Value = 67100
$ws.Range("K83").Value = 65000
$ws.Range("L83").Value = 59900
$ws.Range("M83").Value = 85100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("D89").Value = 707100
$ws.Range("E89").Value = 290000
$ws.Range("F89").Value = 250800
$ws.Range("G89").Value = 352100
$ws.Range("H89").Value = 623600
$ws.Range("I89").Value = 290200
$ws.Range("J89").Value = 101300
$ws.Range("K89").Value = 234500
$ws.Range("L89").Value = 562700
$ws.Range("M89").Value = 103700
$ws.Range("D90").Value = $null
$ws.Range("E90").Value = $null
$ws.Range("F90").Value = $null
$ws.Range("G90").Value = $null
$ws.Range("H90").Value = $null
$ws.Range("I90").Value = $null
$ws.Range("J90").Value = $null
$ws.Range("K90").Value = $null
$ws.Range("L90").Value = $null
$ws.Range("M90").Value = $null
$ws.Range("D91").Value = -87400
$ws.Range("E91").Value = -105300
$ws.Range("F91").Value = -75800
$ws.Range("G91").Value = -60100
$ws.Range("H91").Value = -108800
$ws.Range("I91").Value = -64200
$ws.Range("J91").Value = -51400
$ws.Range("K91").Value = -33300
$ws.Range("L91").Value = -202500
$ws.Range("M91").Value = -232300
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("D94").Value = -530700
$ws.Range("E94").Value = 92000
$ws.Range("F94").Value = -82500
$ws.Range("G94").Value = -981800
$ws.Range("H94").Value = -141500
$ws.Range("I94").Value = -80900
$ws.Range("J94").Value = -65600
$ws.Range("K94").Value = -40700
$ws.Range("L94").Value = -109500
$ws.Range("M94").Value = -81400
$ws.Range("D95").Value = $null
$ws.Range("E95").Value = $null
$ws.Range("F95").Value = $null
$ws.Range("G95").Value = $null
$ws.Range("H95").Value = $null
$ws.Range("I95").Value = $null
$ws.Range("J95").Value = $null
$ws.Range("K95").Value = $null
$ws.Range("L95").Value = $null
$ws.Range("M95").Value = $null
$ws.Range("D96").Value = -147300
$ws.Range("E96").Value = -147300
$ws.Range("F96").Value = -133600
$ws.Range("G96").Value = -134300
$ws.Range("H96").Value = -134400
$ws.Range("I96").Value = -135700
$ws.Range("J96").Value = -128100
$ws.Range("K96").Value = -128000
$ws.Range("L96").Value = -127800
$ws.Range("M96").Value = -128600
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("D100").Value = -412900
$ws.Range("E100").Value = -24100
$ws.Range("F100").Value = -172700
$ws.Range("G100").Value = 725900
$ws.Range("H100").Value = -378000
$ws.Range("I100").Value = -150600
$ws.Range("J100").Value = -58300
$ws.Range("K100").Value = -256900
$ws.Range("L100").Value = -484200
$ws.Range("M100").Value = 90400
$ws.Range("D101").Value = 600
$ws.Range("E101").Value = -1400
$ws.Range("F101").Value = -4600
$ws.Range("G101").Value = 100
$ws.Range("H101").Value = 1100
$ws.Range("I101").Value = 2300
$ws.Range("J101").Value = 1600
$ws.Range("K101").Value = 1200
$ws.Range("L101").Value = -3600
$ws.Range("M101").Value = -1300
$ws.Range("D102").Value = -235800
$ws.Range("E102").Value = 356400
$ws.Range("F102").Value = -9100
$ws.Range("G102").Value = 96300
$ws.Range("H102").Value = 105100
$ws.Range("I102").Value = 61000
$ws.Range("J102").Value = -21000
$ws.Range("K102").Value = -61900
$ws.Range("L102").Value = -36400
$ws.Range("M102").Value = 83100
